# Applies the daily cryptos-list refresh (prices, volumes, and a couple of
# re-ranked coin rows) to Sheet1, cell by cell, matching the GitHub Actions
# scrape commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.212.31'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.33%  '
# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.592.91'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.29%  '
# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.17%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.31'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.11%  '
# Row 6
$ws.Range('E6').Value = '  -0.86%  '
# Row 7
$ws.Range('E7').Value = '  -0.16%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.245'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.12%  '
# Row 9
$ws.Range('E9').Value = '  -0.76%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.03'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.91%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0846'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.16%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.816.14'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.24%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.614.53'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.61%  '
# Row 14
$ws.Range('E14').Value = '  -1.20%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.510'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.94%  '
# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.72'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.95%  '
# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.189.78'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.46%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0726'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.42%  '
# Row 19
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '214.38'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.69%  '
# Row 20
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.36'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.43%  '
# Row 21
$ws.Range('E21').Value = '  -0.07%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.26'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.29%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.05'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.15%  '
# Row 24
$ws.Range('E24').Value = '  -1.60%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.92'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.32%  '
# Row 26
$ws.Range('E26').Value = '  -0.14%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.96'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.20%  '
# Row 28
$ws.Range('E28').Value = '  -1.04%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.10'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.00%  '
# Row 30
$ws.Range('E30').Value = '  -2.27%  '
# Row 31
$ws.Range('E31').Value = '  +0.49%  '
# Row 32
$ws.Range('E32').Value = '  -1.17%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.428.72'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +8.44%  '
# Row 36
$ws.Range('E36').Value = '  -0.60%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.587'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.43%  '
# Row 38
$ws.Range('E38').Value = '  -1.39%  '
# Row 39
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.825'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.48%  '
# Row 40
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.90'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.76%  '
# Row 41
$ws.Range('E41').Value = '  -0.14%  '
# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.982'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -9.49%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.766'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.15%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.13'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.16%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.728.00'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.17%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '61.08'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.94%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '87.05'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.36%  '
# Row 48
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0₆0103'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.12%  '
# Row 49
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.49'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.21%  '
# Row 50
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0502'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.81%  '
# Row 51
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0957'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.93%  '
